$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeiterfassung")

# Copy formatting from row 4 (B=3, C=2, D=4, F=1 styles) to the new rows 7 and 8
$ws.Range("B4:D4").Copy($ws.Range("B7:D7"))
$ws.Range("F4").Copy($ws.Range("F7"))
$ws.Range("B4:D4").Copy($ws.Range("B8:D8"))
$ws.Range("F4").Copy($ws.Range("F8"))

# Row 7 values
$ws.Cells.Item(7, 2).Value = "Pflichtenheft sowie Phasenplan arbeiten"
$ws.Cells.Item(7, 3).Value = "Tobias Lanz"
$ws.Cells.Item(7, 4).Value = 42249
$ws.Cells.Item(7, 6).Value = 2

# Row 8 values
$ws.Cells.Item(8, 2).Value = "Review"
$ws.Cells.Item(8, 3).Value = "Tobias Lanz"
$ws.Cells.Item(8, 4).Value = 42258
$ws.Cells.Item(8, 6).Value = 2

# Update selection to match diff
$ws.Range("O8").Select() | Out-Null
